$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Sfrp1/Fzd2)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sfrp1"
$ws.Range("C2").Value = "Fzd2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.217473
$ws.Range("H2").Value = 0.652419
$ws.Range("I2").Value = 0.003819775075312922
$ws.Range("J2").Value = 0.003819775075312921
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2081856666666667
$ws.Range("N2").Value = 0.624557
$ws.Range("O2").Value = 0.0127208067884984
$ws.Range("P2").Value = 0.0127208067884984
$ws.Range("Q2").Value = 0.04527476148700001
$ws.Range("R2").Value = 0.407472853383
$ws.Range("S2").Value = 0.00004859062070857762
$ws.Range("T2").Value = 0.0000485906207085776

# Row 3: ECs -> FAPs (Sfrp1/Fzd2)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sfrp1"
$ws.Range("C3").Value = "Fzd2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.217473
$ws.Range("H3").Value = 0.652419
$ws.Range("I3").Value = 0.003819775075312922
$ws.Range("J3").Value = 0.003819775075312921
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 13.42533133333333
$ws.Range("N3").Value = 40.275994
$ws.Range("O3").Value = 0.8203304708596988
$ws.Range("P3").Value = 0.8203304708596987
$ws.Range("Q3").Value = 2.919647081054
$ws.Range("R3").Value = 26.276823729486
$ws.Range("S3").Value = 0.003133477886109591
$ws.Range("T3").Value = 0.00313347788610959

# Row 4: ECs -> M2 (Sfrp1/Fzd2)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sfrp1"
$ws.Range("C4").Value = "Fzd2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.217473
$ws.Range("H4").Value = 0.652419
$ws.Range("I4").Value = 0.003819775075312922
$ws.Range("J4").Value = 0.003819775075312921
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04883999999999999
$ws.Range("N4").Value = 0.14652
$ws.Range("O4").Value = 0.002984279434304292
$ws.Range("P4").Value = 0.002984279434304292
$ws.Range("Q4").Value = 0.01062138132
$ws.Range("R4").Value = 0.09559243187999998
$ws.Range("S4").Value = 0.00001139927620092448
$ws.Range("T4").Value = 0.00001139927620092448

# Row 5: ECs -> sCs (Sfrp1/Fzd2)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sfrp1"
$ws.Range("C5").Value = "Fzd2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.217473
$ws.Range("H5").Value = 0.652419
$ws.Range("I5").Value = 0.003819775075312922
$ws.Range("J5").Value = 0.003819775075312921
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.683402666666666
$ws.Range("N5").Value = 8.050208
$ws.Range("O5").Value = 0.1639644429174985
$ws.Range("P5").Value = 0.1639644429174985
$ws.Range("Q5").Value = 0.5835676281279999
$ws.Range("R5").Value = 5.252108653152
$ws.Range("S5").Value = 0.0006263072922938293
$ws.Range("T5").Value = 0.0006263072922938292

# Row 6: FAPs -> ECs (Sfrp1/Fzd2)
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sfrp1"
$ws.Range("C6").Value = "Fzd2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 55.43187966666667
$ws.Range("H6").Value = 166.295639
$ws.Range("I6").Value = 0.9736257481548445
$ws.Range("J6").Value = 0.9736257481548444
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2081856666666667
$ws.Range("N6").Value = 0.624557
$ws.Range("O6").Value = 0.0127208067884984
$ws.Range("P6").Value = 0.0127208067884984
$ws.Range("Q6").Value = 11.54012282299145
$ws.Range("R6").Value = 103.861105406923
$ws.Range("S6").Value = 0.01238530502658498
$ws.Range("T6").Value = 0.01238530502658498

# Row 7: FAPs -> FAPs (Sfrp1/Fzd2)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sfrp1"
$ws.Range("C7").Value = "Fzd2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 55.43187966666667
$ws.Range("H7").Value = 166.295639
$ws.Range("I7").Value = 0.9736257481548445
$ws.Range("J7").Value = 0.9736257481548444
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.42533133333333
$ws.Range("N7").Value = 40.275994
$ws.Range("O7").Value = 0.8203304708596988
$ws.Range("P7").Value = 0.8203304708596987
$ws.Range("Q7").Value = 744.1913509544628
$ws.Range("R7").Value = 6697.722158590165
$ws.Range("S7").Value = 0.79869486842499
$ws.Range("T7").Value = 0.7986948684249899

# Row 8: FAPs -> M2 (Sfrp1/Fzd2)
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sfrp1"
$ws.Range("C8").Value = "Fzd2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 55.43187966666667
$ws.Range("H8").Value = 166.295639
$ws.Range("I8").Value = 0.9736257481548445
$ws.Range("J8").Value = 0.9736257481548444
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.04883999999999999
$ws.Range("N8").Value = 0.14652
$ws.Range("O8").Value = 0.002984279434304292
$ws.Range("P8").Value = 0.002984279434304292
$ws.Range("Q8").Value = 2.70729300292
$ws.Range("R8").Value = 24.36563702628
$ws.Range("S8").Value = 0.002905571296927632
$ws.Range("T8").Value = 0.002905571296927632

# Row 9: FAPs -> sCs (Sfrp1/Fzd2)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sfrp1"
$ws.Range("C9").Value = "Fzd2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 55.43187966666667
$ws.Range("H9").Value = 166.295639
$ws.Range("I9").Value = 0.9736257481548445
$ws.Range("J9").Value = 0.9736257481548444
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.683402666666666
$ws.Range("N9").Value = 8.050208
$ws.Range("O9").Value = 0.1639644429174985
$ws.Range("P9").Value = 0.1639644429174985
$ws.Range("Q9").Value = 148.7460537158791
$ws.Range("R9").Value = 1338.714483442912
$ws.Range("S9").Value = 0.1596400034063418
$ws.Range("T9").Value = 0.1596400034063418

# Row 10: sCs -> ECs (Sfrp1/Fzd2)
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sfrp1"
$ws.Range("C10").Value = "Fzd2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.284104333333334
$ws.Range("H10").Value = 3.852313000000001
$ws.Range("I10").Value = 0.02255447676984262
$ws.Range("J10").Value = 0.02255447676984261
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2081856666666667
$ws.Range("N10").Value = 0.624557
$ws.Range("O10").Value = 0.0127208067884984
$ws.Range("P10").Value = 0.0127208067884984
$ws.Range("Q10").Value = 0.2673321167045556
$ws.Range("R10").Value = 2.405989050341001
$ws.Range("S10").Value = 0.0002869111412048435
$ws.Range("T10").Value = 0.0002869111412048434

# Row 11: sCs -> FAPs (Sfrp1/Fzd2)
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Sfrp1"
$ws.Range("C11").Value = "Fzd2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.284104333333334
$ws.Range("H11").Value = 3.852313000000001
$ws.Range("I11").Value = 0.02255447676984262
$ws.Range("J11").Value = 0.02255447676984261
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 13.42533133333333
$ws.Range("N11").Value = 40.275994
$ws.Range("O11").Value = 0.8203304708596988
$ws.Range("P11").Value = 0.8203304708596987
$ws.Range("Q11").Value = 17.23952614156911
$ws.Range("R11").Value = 155.155735274122
$ws.Range("S11").Value = 0.01850212454859913
$ws.Range("T11").Value = 0.01850212454859913

# Row 12: sCs -> M2 (Sfrp1/Fzd2)
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Sfrp1"
$ws.Range("C12").Value = "Fzd2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.284104333333334
$ws.Range("H12").Value = 3.852313000000001
$ws.Range("I12").Value = 0.02255447676984262
$ws.Range("J12").Value = 0.02255447676984261
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.04883999999999999
$ws.Range("N12").Value = 0.14652
$ws.Range("O12").Value = 0.002984279434304292
$ws.Range("P12").Value = 0.002984279434304292
$ws.Range("Q12").Value = 0.06271565564000001
$ws.Range("R12").Value = 0.56444090076
$ws.Range("S12").Value = 0.00006730886117573522
$ws.Range("T12").Value = 0.0000673088611757352

# Row 13: sCs -> sCs (Sfrp1/Fzd2)
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Sfrp1"
$ws.Range("C13").Value = "Fzd2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.284104333333334
$ws.Range("H13").Value = 3.852313000000001
$ws.Range("I13").Value = 0.02255447676984262
$ws.Range("J13").Value = 0.02255447676984261
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.683402666666666
$ws.Range("N13").Value = 8.050208
$ws.Range("O13").Value = 0.1639644429174985
$ws.Range("P13").Value = 0.1639644429174985
$ws.Range("Q13").Value = 3.445768992344889
$ws.Range("R13").Value = 31.011920931104
$ws.Range("S13").Value = 0.003698132218862907
$ws.Range("T13").Value = 0.003698132218862906
